$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in column H, row 1, matching the formatting of the
# other header cells (e.g. G1) so it reuses the existing bold/bordered style.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Add the corresponding "Save" value for the data row.
$ws.Range("H2").Value = 1
